# Generate Report for Handback
# Updates the Status of the "01ca06e2-a6df-4437-bb8d-b2f1c508ac42.md" file
# from "Ready for handoff" to "Handback transform failed" across the
# Overview, zh-cn, and de-de sheets, and records the transform-failure
# Error Detail message for that file on the zh-cn and de-de sheets.

$wb = $excel.ActiveWorkbook

# --- Overview sheet ---
$wsOverview = $wb.Worksheets.Item("Overview")
$wsOverview.Range("B3").Value = "Handback transform failed"
$wsOverview.Range("C3").Value = "Handback transform failed"

# --- zh-cn sheet ---
$wsZhCn = $wb.Worksheets.Item("zh-cn")
$wsZhCn.Range("C3").Value = "Handback transform failed"
$wsZhCn.Range("K3").Value = "Handback file name: xhboafbd.wp5 is different with handoff file name: 01ca06e2-a6df-4437-bb8d-b2f1c508ac42.51a522d05f75cd64ac1a0c14d51443e3d822e4d5.zh-cn."

# --- de-de sheet ---
$wsDeDe = $wb.Worksheets.Item("de-de")
$wsDeDe.Range("C3").Value = "Handback transform failed"
$wsDeDe.Range("K3").Value = "Handback file name: xhboafbd.wp5 is different with handoff file name: 01ca06e2-a6df-4437-bb8d-b2f1c508ac42.51a522d05f75cd64ac1a0c14d51443e3d822e4d5.de-de."
